$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) per row, matching the diff.
# Row 8, 24, 25 are intentionally left unchanged.
$updates = @{
    2  = 4
    3  = -6
    4  = -7
    5  = -3
    6  = -2
    7  = -5
    9  = -3
    10 = -1
    11 = -1
    12 = 5
    13 = -6
    14 = 2
    15 = 3
    16 = -5
    17 = -2
    18 = 2
    19 = -4
    20 = -1
    21 = 1
    22 = -2
    23 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
